$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 6384
$ws.Range("L3").Value = 6885
$ws.Range("L4").Value = 1706
$ws.Range("L5").Value = 405
$ws.Range("L6").Value = 5656
$ws.Range("L7").Value = 21036

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L6").Value = 169
$ws.Range("L7").Value = 669
$ws.Range("L8").Value = 1391
$ws.Range("L11").Value = 348
$ws.Range("L12").Value = 48
$ws.Range("L17").Value = 37
$ws.Range("L19").Value = 572
$ws.Range("L29").Value = 1176
$ws.Range("L33").Value = 948
$ws.Range("L36").Value = 268
$ws.Range("L37").Value = 807
$ws.Range("L42").Value = 664
$ws.Range("L46").Value = 51
$ws.Range("L48").Value = 274
$ws.Range("L51").Value = 262
$ws.Range("L52").Value = 449
$ws.Range("L59").Value = 35
$ws.Range("L64").Value = 134
$ws.Range("L65").Value = 414
$ws.Range("L67").Value = 728
$ws.Range("L68").Value = 66
$ws.Range("L77").Value = 140
$ws.Range("L78").Value = 278
$ws.Range("L79").Value = 583
$ws.Range("L83").Value = 460
$ws.Range("L88").Value = 223
$ws.Range("L91").Value = 282
$ws.Range("L92").Value = 65
$ws.Range("L94").Value = 256
$ws.Range("L101").Value = 21036

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 169

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L3").Value = 214
$ws.Range("L7").Value = 669

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 425
$ws.Range("L3").Value = 492
$ws.Range("L6").Value = 336
$ws.Range("L7").Value = 1391

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 131
$ws.Range("L7").Value = 348

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("L2").Value = 206
$ws.Range("L3").Value = 175
$ws.Range("L7").Value = 572

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("L3").Value = 453
$ws.Range("L6").Value = 286
$ws.Range("L7").Value = 1176

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("L2").Value = 261
$ws.Range("L4").Value = 64
$ws.Range("L6").Value = 267
$ws.Range("L7").Value = 948

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L3").Value = 88
$ws.Range("L7").Value = 268

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L2").Value = 246
$ws.Range("L7").Value = 807

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L6").Value = 189
$ws.Range("L7").Value = 664

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("L3").Value = 14
$ws.Range("L7").Value = 51

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 71
$ws.Range("L7").Value = 274

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("L2").Value = 81
$ws.Range("L4").Value = 38
$ws.Range("L5").Value = 7
$ws.Range("L7").Value = 262

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L2").Value = 143
$ws.Range("L7").Value = 449

$ws = $wb.Worksheets.Item('Montclare')
$ws.Range("L3").Value = 15
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 151
$ws.Range("L7").Value = 414

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L3").Value = 283
$ws.Range("L7").Value = 728

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 66

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("L3").Value = 44
$ws.Range("L7").Value = 140

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("L2").Value = 72
$ws.Range("L6").Value = 79
$ws.Range("L7").Value = 278

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L3").Value = 188
$ws.Range("L7").Value = 583

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L2").Value = 148
$ws.Range("L7").Value = 460

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("L6").Value = 61
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("L3").Value = 127
$ws.Range("L7").Value = 282

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L3").Value = 12
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("L2").Value = 62
$ws.Range("L7").Value = 256
